{"js": "// Append a new \"To Do List\" bullet item at the very end of the document,\n// right after the paragraph ending in \"...possible to add parallax backgrounds\".\n// The new paragraph keeps the same list (numId 2 / ListParagraph style) as the\n// other \"To do List\" bullets already in the document, matching their run\n// formatting (sz/szCs 24).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// insertParagraph inherits the paragraph/list formatting of lastParagraph,\n// which is exactly the \"To do List\" bullet (pStyle ListParagraph, numId 2).\nlastParagraph.insertParagraph(\n  \"Index load screen just like circle language button (Load circle) with a spinning light in its orbit\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Append a new \"To Do List\" bullet item at the very end of the document,\n# right after the paragraph ending in \"...possible to add parallax backgrounds\".\n# The new paragraph is inserted via InsertParagraphAfter on the last\n# paragraph's range, so it automatically inherits the same list/paragraph\n# formatting (pStyle ListParagraph, numId 2) and run formatting (sz/szCs 24)\n# as the existing \"To do List\" bullets.\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastRange = $lastParagraph.Range\n$lastRange.Collapse(0)  # wdCollapseEnd\n$lastRange.InsertParagraphAfter()\n\n$newRange = $d.Paragraphs.Last.Range\n$newRange.InsertAfter(\"Index load screen just like circle language button (Load circle) with a spinning light in its orbit\")\n"}
